# Remove the two introductory "about TED" paragraphs from the content
# placeholder on slide 4, leaving only the "TED² is a navigational tool..."
# paragraph (commit: "removal of info on presentation").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# The first paragraph ("TED - Technology, Entertainment and Design. ") and
# the second paragraph ("TED talks are a collection of brief yet ...") are
# deleted in place; after each delete the remaining text shifts up so the
# target paragraph is always back at index 1.
$tr.Paragraphs(1, 1).Delete()
$tr.Paragraphs(1, 1).Delete()
